$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.445.71'
$ws.Range("E2").Value = '  +2.71%  '
$ws.Range("D3").Value = '1.676.04'
$ws.Range("E3").Value = '  +4.06%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5309'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.69%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2672'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06394'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.55'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07807'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.86%  '
$ws.Range("D12").Value = '1.673.49'
$ws.Range("E12").Value = '  +3.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.499'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5562'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.94%  '
$ws.Range("D15").Value = '0.0₅8347'
$ws.Range("E15").Value = '  +6.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.64'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.25%  '
$ws.Range("D17").Value = '26.483.33'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.002'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.773'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '195.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.88%  '
$ws.Range("E21").Value = '  +3.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.337'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '144.05'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1277'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.424'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.426'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06143'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.273'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.616'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.448'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.691'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.37%  '
$ws.Range("E34").Value = '  +4.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.425'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.84%  '
$ws.Range("E36").Value = '  +2.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5721'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01638'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.033'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.00%  '
$ws.Range("D40").Value = '1.071.93'
$ws.Range("E40").Value = '  +5.66%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8584'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.62%  '
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.02'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.17%  '
$ws.Range("D44").Value = '1.825.08'
$ws.Range("E44").Value = '  +3.80%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.09'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.85%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.0₈107'
$ws.Range("E46").Value = '  -2.63%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.163'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.80%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.006'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05204'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.05%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.472'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.75%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.027'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.76%  '
